$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H2").Value = 202.42857
$ws.Range("I2").Value = 207
$ws.Range("J2").Value = 199.88889
$ws.Range("K2").Value = 207
$ws.Range("L2").Value = 199.88889
$ws.Range("M2").Value = -94
$ws.Range("N2").Value = -425.88889
$ws.Range("H6").Value = 1
$ws.Range("I6").Value = 1
$ws.Range("J6").Value = 0
$ws.Range("K6").Value = 3
$ws.Range("L6").Value = 0
$ws.Range("M6").Value = 109
$ws.Range("N6").ClearContents()
$ws.Range("H9").Value = 59499.35
$ws.Range("I9").Value = 91058.37
$ws.Range("J9").Value = 1641.1666
$ws.Range("K9").Value = 91058.37
$ws.Range("L9").Value = 1641.1666
$ws.Range("M9").Value = -90889.37
$ws.Range("N9").Value = -1979.1666
$ws.Range("H17").Value = 633.4400000000001
$ws.Range("J17").Value = 633.4400000000001
$ws.Range("L17").Value = 1900.32
$ws.Range("N17").Value = -2236.32
$ws.Range("H19").Value = 768.8
$ws.Range("J19").Value = 786
$ws.Range("L19").Value = 786
$ws.Range("N19").Value = -1136
$ws.Range("H38").Value = 332.6
$ws.Range("I38").Value = 41.25
$ws.Range("J38").Value = 1498
$ws.Range("K38").Value = 123.75
$ws.Range("L38").Value = 4494
$ws.Range("M38").Value = 248.25
$ws.Range("N38").Value = -5238
$ws.Range("H40").Value = 2212.3125
$ws.Range("I40").Value = 2699
$ws.Range("J40").Value = 2142.7856
$ws.Range("K40").Value = 2699
$ws.Range("L40").Value = 2142.7856
$ws.Range("M40").Value = -2524
$ws.Range("N40").Value = -2492.7856
$ws.Range("H74").Value = 6359.2593
$ws.Range("I74").Value = 4970
$ws.Range("K74").Value = 4970
$ws.Range("M74").Value = -4034
$ws.Range("H75").Value = 117416.5
$ws.Range("J75").Value = 117416.5
$ws.Range("L75").Value = 117416.5
$ws.Range("N75").Value = -119288.5
$ws.Range("H76").Value = 72146264
$ws.Range("J76").Value = 4467.857
$ws.Range("L76").Value = 4467.857
$ws.Range("N76").Value = -5097.857
$ws.Range("H77").Value = 6359.2593
$ws.Range("I77").Value = 4970
$ws.Range("K77").Value = 24850
$ws.Range("M77").Value = -20170
$ws.Range("H78").Value = 117416.5
$ws.Range("J78").Value = 117416.5
$ws.Range("L78").Value = 352249.5
$ws.Range("N78").Value = -361609.5
$ws.Range("H79").Value = 72146264
$ws.Range("J79").Value = 4467.857
$ws.Range("L79").Value = 4467.857
$ws.Range("N79").Value = -6651.857
$ws.Range("H80").Value = 927.44446
$ws.Range("I80").Value = 999.4
$ws.Range("K80").Value = 2998.2
$ws.Range("M80").Value = -2000.2
$ws.Range("H83").Value = 927.44446
$ws.Range("I83").Value = 999.4
$ws.Range("K83").Value = 8994.6
$ws.Range("M83").Value = -4002.6
$ws.Range("H98").Value = 0
$ws.Range("I98").Value = 0
$ws.Range("J98").Value = 0
$ws.Range("K98").Value = 0
$ws.Range("L98").Value = 0
$ws.Range("M98").ClearContents()
$ws.Range("N98").ClearContents()
$ws.Range("H100").Value = 3968.125
$ws.Range("I100").Value = 1686.375
$ws.Range("K100").Value = 1686.375
$ws.Range("M100").Value = -1145.375
$ws.Range("H107").Value = 2339.36
$ws.Range("I107").Value = 1563.6522
$ws.Range("J107").Value = 11260
$ws.Range("K107").Value = 1563.6522
$ws.Range("L107").Value = 11260
$ws.Range("M107").Value = 356.3478
$ws.Range("N107").Value = -15100
$ws.Range("H111").Value = 2939
$ws.Range("I111").Value = 2433.7646
$ws.Range("K111").Value = 7301.293799999999
$ws.Range("M111").Value = -4234.293799999999
$ws.Range("H122").Value = 0
$ws.Range("I122").Value = 0
$ws.Range("J122").Value = 0
$ws.Range("K122").Value = 0
$ws.Range("L122").Value = 0
$ws.Range("M122").ClearContents()
$ws.Range("N122").ClearContents()
$ws.Range("H134").Value = 43726.273
$ws.Range("J134").Value = 43726.273
$ws.Range("L134").Value = 43726.273
$ws.Range("N134").Value = -53866.273
$ws.Range("H137").Value = 2136.5
$ws.Range("I137").Value = 1340.6666
$ws.Range("J137").Value = 2733.375
$ws.Range("K137").Value = 4021.9998
$ws.Range("L137").Value = 8200.125
$ws.Range("M137").Value = -1471.9998
$ws.Range("N137").Value = -13300.125
$ws.Range("H138").Value = 3802.543
$ws.Range("I138").Value = 3240.4285
$ws.Range("K138").Value = 9721.2855
$ws.Range("M138").Value = -4581.2855
$ws.Range("H141").Value = 7052.0454
$ws.Range("I141").Value = 4120.0586
$ws.Range("K141").Value = 12360.1758
$ws.Range("M141").Value = -7180.175800000001
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 3157.7556
$ws.Range("I32").Value = 2740.175
$ws.Range("K32").Value = 2740.175
$ws.Range("M32").Value = -2453.175
$ws.Range("H45").Value = 2070.1428
$ws.Range("I45").Value = 1436.25
$ws.Range("K45").Value = 1436.25
$ws.Range("M45").Value = -1059.25
$ws.Range("H61").Value = 13113.056
$ws.Range("I61").Value = 10677.5
$ws.Range("J61").Value = 16157.5
$ws.Range("K61").Value = 10677.5
$ws.Range("L61").Value = 16157.5
$ws.Range("M61").Value = -10465.5
$ws.Range("N61").Value = -16581.5
$ws.Range("H74").Value = 4713.8667
$ws.Range("I74").Value = 4269.909
$ws.Range("K74").Value = 4269.909
$ws.Range("M74").Value = -3395.909
$ws.Range("H77").Value = 4713.8667
$ws.Range("I77").Value = 4269.909
$ws.Range("K77").Value = 21349.545
$ws.Range("M77").Value = -16981.545
$ws.Range("H88").Value = 4897.067
$ws.Range("I88").Value = 1129.8889
$ws.Range("J88").Value = 6511.5713
$ws.Range("K88").Value = 1129.8889
$ws.Range("L88").Value = 6511.5713
$ws.Range("M88").Value = -723.8888999999999
$ws.Range("N88").Value = -7323.5713
$ws.Range("H91").Value = 4897.067
$ws.Range("I91").Value = 1129.8889
$ws.Range("J91").Value = 6511.5713
$ws.Range("K91").Value = 1129.8889
$ws.Range("L91").Value = 6511.5713
$ws.Range("M91").Value = 274.1111000000001
$ws.Range("N91").Value = -9319.5713
$ws.Range("H92").Value = 2683.3333
$ws.Range("J92").Value = 2683.3333
$ws.Range("L92").Value = 2683.3333
$ws.Range("N92").Value = -7675.3333
$ws.Range("H102").Value = 2190.8
$ws.Range("I102").Value = 2190.8
$ws.Range("K102").Value = 2190.8
$ws.Range("M102").Value = -568.8000000000002
$ws.Range("H110").Value = 1687.6666
$ws.Range("I110").Value = 1594.5
$ws.Range("J110").Value = 1874
$ws.Range("K110").Value = 1594.5
$ws.Range("L110").Value = 1874
$ws.Range("M110").Value = 450.5
$ws.Range("N110").Value = -5964
$ws.Range("H132").Value = 3616.724
$ws.Range("I132").Value = 2333.8635
$ws.Range("K132").Value = 7001.5905
$ws.Range("M132").Value = -4471.5905
$ws.Range("H136").Value = 13113.056
$ws.Range("I136").Value = 10677.5
$ws.Range("J136").Value = 16157.5
$ws.Range("K136").Value = 32032.5
$ws.Range("L136").Value = 48472.5
$ws.Range("M136").Value = -29482.5
$ws.Range("N136").Value = -53572.5
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H75").Value = 11000
$ws.Range("I75").Value = 11000
$ws.Range("J75").Value = 0
$ws.Range("K75").Value = 11000
$ws.Range("L75").Value = 0
$ws.Range("M75").Value = -10064
$ws.Range("N75").ClearContents()
$ws.Range("H78").Value = 11000
$ws.Range("I78").Value = 11000
$ws.Range("J78").Value = 0
$ws.Range("K78").Value = 33000
$ws.Range("L78").Value = 0
$ws.Range("M78").Value = -28320
$ws.Range("N78").ClearContents()
$ws.Range("H86").Value = 10392.818
$ws.Range("I86").Value = 1103
$ws.Range("J86").Value = 35165.668
$ws.Range("K86").Value = 1103
$ws.Range("L86").Value = 35165.668
$ws.Range("M86").Value = 20
$ws.Range("N86").Value = -37411.668
$ws.Range("H89").Value = 10392.818
$ws.Range("I89").Value = 1103
$ws.Range("J89").Value = 35165.668
$ws.Range("K89").Value = 5515
$ws.Range("L89").Value = 175828.34
$ws.Range("M89").Value = 101
$ws.Range("N89").Value = -187060.34
$ws.Range("H94").Value = 7159405
$ws.Range("I94").Value = 20001516
$ws.Range("K94").Value = 20001516
$ws.Range("M94").Value = -20001065
$ws.Range("H105").Value = 2398.3157
$ws.Range("I105").Value = 2113
$ws.Range("J105").Value = 3016.5
$ws.Range("K105").Value = 2113
$ws.Range("L105").Value = 3016.5
$ws.Range("M105").Value = -366
$ws.Range("N105").Value = -6510.5
$ws.Range("H107").Value = 1671.1111
$ws.Range("I107").Value = 1732.2106
$ws.Range("J107").Value = 1339.4286
$ws.Range("K107").Value = 1732.2106
$ws.Range("L107").Value = 1339.4286
$ws.Range("M107").Value = 187.7893999999999
$ws.Range("N107").Value = -5179.4286
$ws.Range("H132").Value = 64363
$ws.Range("J132").Value = 67235.60000000001
$ws.Range("L132").Value = 67235.60000000001
$ws.Range("N132").Value = -77355.60000000001
$ws.Range("H133").Value = 0
$ws.Range("J133").Value = 0
$ws.Range("L133").Value = 0
$ws.Range("N133").ClearContents()
$ws.Range("H134").Value = 11387.173
$ws.Range("I134").Value = 11973.25
$ws.Range("J134").Value = 10665.846
$ws.Range("K134").Value = 35919.75
$ws.Range("L134").Value = 31997.538
$ws.Range("M134").Value = -33384.75
$ws.Range("N134").Value = -37067.538
$ws.Range("H139").Value = 208445
$ws.Range("J139").Value = 208445
$ws.Range("L139").Value = 208445
$ws.Range("N139").Value = -218725
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H16").Value = 1233.4
$ws.Range("J16").Value = 1888
$ws.Range("L16").Value = 1888
$ws.Range("N16").Value = -2462
$ws.Range("H22").Value = 389.12195
$ws.Range("J22").Value = 558.9474
$ws.Range("L22").Value = 558.9474
$ws.Range("N22").Value = -1258.9474
$ws.Range("H31").Value = 3092.92
$ws.Range("I31").Value = 1322.5555
$ws.Range("J31").Value = 4088.75
$ws.Range("K31").Value = 1322.5555
$ws.Range("L31").Value = 4088.75
$ws.Range("M31").Value = -1027.5555
$ws.Range("N31").Value = -4678.75
$ws.Range("H34").Value = 3092.92
$ws.Range("I34").Value = 1322.5555
$ws.Range("J34").Value = 4088.75
$ws.Range("K34").Value = 1322.5555
$ws.Range("L34").Value = 4088.75
$ws.Range("M34").Value = -1120.5555
$ws.Range("N34").Value = -4492.75
$ws.Range("H35").Value = 3666.3333
$ws.Range("I35").Value = 3666.3333
$ws.Range("K35").Value = 3666.3333
$ws.Range("M35").Value = -3372.3333
$ws.Range("H41").Value = 15496.5
$ws.Range("J41").Value = 25325
$ws.Range("L41").Value = 25325
$ws.Range("N41").Value = -26181
$ws.Range("H60").Value = 32333
$ws.Range("J60").Value = 34999.617
$ws.Range("L60").Value = 34999.617
$ws.Range("N60").Value = -36021.617
$ws.Range("H94").Value = 2227.0715
$ws.Range("I94").Value = 2118.3333
$ws.Range("J94").Value = 2308.625
$ws.Range("K94").Value = 2118.3333
$ws.Range("L94").Value = 2308.625
$ws.Range("M94").Value = -1667.3333
$ws.Range("N94").Value = -3210.625
$ws.Range("H105").Value = 1465.84
$ws.Range("I105").Value = 1254.6666
$ws.Range("K105").Value = 1254.6666
$ws.Range("M105").Value = 492.3334
$ws.Range("H113").Value = 1233.4
$ws.Range("J113").Value = 1888
$ws.Range("L113").Value = 1888
$ws.Range("N113").Value = -6228
$ws.Range("H134").Value = 4923.683
$ws.Range("I134").Value = 4319
$ws.Range("K134").Value = 12957
$ws.Range("M134").Value = -10422
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H39").Value = 4282.9165
$ws.Range("I39").Value = 325
$ws.Range("J39").Value = 4455
$ws.Range("K39").Value = 975
$ws.Range("L39").Value = 13365
$ws.Range("M39").Value = -681
$ws.Range("N39").Value = -13953
$ws.Range("H45").Value = 331
$ws.Range("J45").Value = 333
$ws.Range("L45").Value = 999
$ws.Range("N45").Value = -2063
$ws.Range("H62").Value = 8349.700000000001
$ws.Range("J62").Value = 8874.625
$ws.Range("L62").Value = 26623.875
$ws.Range("N62").Value = -27995.875
$ws.Range("H65").Value = 8349.700000000001
$ws.Range("J65").Value = 8874.625
$ws.Range("L65").Value = 79871.625
$ws.Range("N65").Value = -86735.625
$ws.Range("H75").Value = 4665.3335
$ws.Range("J75").Value = 3999
$ws.Range("L75").Value = 11997
$ws.Range("N75").Value = -13993
$ws.Range("H78").Value = 4665.3335
$ws.Range("J78").Value = 3999
$ws.Range("L78").Value = 35991
$ws.Range("N78").Value = -45975
$ws.Range("H98").Value = 38466196
$ws.Range("J98").Value = 50004708
$ws.Range("L98").Value = 150014124
$ws.Range("N98").Value = -150017120
$ws.Range("H120").Value = 15404.75
$ws.Range("I120").Value = 13319.857
$ws.Range("K120").Value = 39959.571
$ws.Range("M120").Value = -35121.571
$ws.Range("H122").Value = 499.45
$ws.Range("J122").Value = 447.5
$ws.Range("L122").Value = 4027.5
$ws.Range("N122").Value = -8927.5
$ws.Range("H131").Value = 3596.3855
$ws.Range("I131").Value = 2350
$ws.Range("J131").Value = 3650.5762
$ws.Range("K131").Value = 7050
$ws.Range("L131").Value = 10951.7286
$ws.Range("M131").Value = -2010
$ws.Range("N131").Value = -21031.7286
$ws.Range("H140").Value = 8934618
$ws.Range("I140").Value = 14708529
$ws.Range("K140").Value = 44125587
$ws.Range("M140").Value = -44120407
$ws.Range("H141").Value = 28400.5
$ws.Range("I141").Value = 4908.923
$ws.Range("K141").Value = 14726.769
$ws.Range("M141").Value = -9546.769
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H59").Value = 0
$ws.Range("J59").Value = 0
$ws.Range("L59").Value = 0
$ws.Range("N59").ClearContents()
$ws.Range("H80").Value = 2163.2856
$ws.Range("I80").Value = 2335.75
$ws.Range("J80").Value = 1933.3334
$ws.Range("K80").Value = 2335.75
$ws.Range("L80").Value = 1933.3334
$ws.Range("M80").Value = -1337.75
$ws.Range("N80").Value = -3929.3334
$ws.Range("H82").Value = 75000
$ws.Range("I82").Value = 75000
$ws.Range("K82").Value = 75000
$ws.Range("M82").Value = -74617
$ws.Range("H83").Value = 2163.2856
$ws.Range("I83").Value = 2335.75
$ws.Range("J83").Value = 1933.3334
$ws.Range("K83").Value = 11678.75
$ws.Range("L83").Value = 9666.666999999999
$ws.Range("M83").Value = -6686.75
$ws.Range("N83").Value = -19650.667
$ws.Range("H85").Value = 75000
$ws.Range("I85").Value = 75000
$ws.Range("K85").Value = 75000
$ws.Range("M85").Value = -73674
$ws.Range("H98").Value = 17194.4
$ws.Range("J98").Value = 17194.4
$ws.Range("L98").Value = 17194.4
$ws.Range("N98").Value = -23184.4
$ws.Range("H107").Value = 262.5
$ws.Range("I107").Value = 190.5
$ws.Range("J107").Value = 370.5
$ws.Range("K107").Value = 190.5
$ws.Range("L107").Value = 370.5
$ws.Range("M107").Value = 1729.5
$ws.Range("N107").Value = -4210.5
$ws.Range("H132").Value = 5405.3335
$ws.Range("I132").Value = 3552.5715
$ws.Range("K132").Value = 10657.7145
$ws.Range("M132").Value = -8127.7145
$ws.Range("H135").Value = 40000
$ws.Range("I135").Value = 40000
$ws.Range("K135").Value = 40000
$ws.Range("M135").Value = -34930
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H7").Value = 501004
$ws.Range("J7").Value = 0
$ws.Range("L7").Value = 0
$ws.Range("N7").ClearContents()
$ws.Range("H16").Value = 0
$ws.Range("J16").Value = 0
$ws.Range("L16").Value = 0
$ws.Range("N16").ClearContents()
$ws.Range("H22").Value = 994.2632
$ws.Range("J22").Value = 1058.3125
$ws.Range("L22").Value = 1058.3125
$ws.Range("N22").Value = -1648.3125
$ws.Range("H27").Value = 994.2632
$ws.Range("J27").Value = 1058.3125
$ws.Range("L27").Value = 1058.3125
$ws.Range("N27").Value = -1272.3125
$ws.Range("H46").Value = 2106.3333
$ws.Range("J46").Value = 2460.6667
$ws.Range("L46").Value = 2460.6667
$ws.Range("N46").Value = -2836.6667
$ws.Range("H55").Value = 414.92856
$ws.Range("I55").Value = 140.44444
$ws.Range("K55").Value = 140.44444
$ws.Range("M55").Value = 32.55556000000001
$ws.Range("H57").Value = 24999
$ws.Range("J57").Value = 29998
$ws.Range("L57").Value = 29998
$ws.Range("N57").Value = -31130
$ws.Range("H62").Value = 23333.666
$ws.Range("J62").Value = 30000
$ws.Range("L62").Value = 30000
$ws.Range("N62").Value = -31248
$ws.Range("H65").Value = 23333.666
$ws.Range("J65").Value = 30000
$ws.Range("L65").Value = 90000
$ws.Range("N65").Value = -96240
$ws.Range("H100").Value = 1592.2
$ws.Range("I100").Value = 1203.1428
$ws.Range("K100").Value = 1203.1428
$ws.Range("M100").Value = -662.1428000000001
$ws.Range("H126").Value = 501004
$ws.Range("J126").Value = 0
$ws.Range("L126").Value = 0
$ws.Range("N126").ClearContents()
$ws.Range("H130").Value = 63143
$ws.Range("J130").Value = 54714.5
$ws.Range("L130").Value = 54714.5
$ws.Range("N130").Value = -64754.5
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H21").Value = 0
$ws.Range("J21").Value = 0
$ws.Range("L21").Value = 0
$ws.Range("N21").ClearContents()
$ws.Range("H31").Value = 0
$ws.Range("J31").Value = 0
$ws.Range("L31").Value = 0
$ws.Range("N31").ClearContents()
$ws.Range("H35").Value = 0
$ws.Range("J35").Value = 0
$ws.Range("L35").Value = 0
$ws.Range("N35").ClearContents()
$ws.Range("H69").Value = 37166.668
$ws.Range("J69").Value = 37166.668
$ws.Range("L69").Value = 37166.668
$ws.Range("N69").Value = -38664.668
$ws.Range("H72").Value = 37166.668
$ws.Range("J72").Value = 37166.668
$ws.Range("L72").Value = 111500.004
$ws.Range("N72").Value = -118988.004
$ws.Range("H80").Value = 80000.336
$ws.Range("J80").Value = 80000.336
$ws.Range("L80").Value = 80000.336
$ws.Range("N80").Value = -81996.336
$ws.Range("H83").Value = 80000.336
$ws.Range("J83").Value = 80000.336
$ws.Range("L83").Value = 240001.008
$ws.Range("N83").Value = -249985.008
$ws.Range("H92").Value = 49999
$ws.Range("J92").Value = 49999
$ws.Range("L92").Value = 49999
$ws.Range("N92").Value = -54991
$ws.Range("H126").Value = 2153.182
$ws.Range("I126").Value = 2185.625
$ws.Range("J126").Value = 2066.6667
$ws.Range("K126").Value = 6556.875
$ws.Range("L126").Value = 6200.000100000001
$ws.Range("M126").Value = -4086.875
$ws.Range("N126").Value = -11140.0001
$ws.Range("H132").Value = 6028.7617
$ws.Range("I132").Value = 6350.5
$ws.Range("K132").Value = 19051.5
$ws.Range("M132").Value = -16521.5
$ws.Range("H136").Value = 16646.066
$ws.Range("I136").Value = 16646.066
$ws.Range("J136").Value = 0
$ws.Range("K136").Value = 49938.198
$ws.Range("L136").Value = 0
$ws.Range("M136").Value = -47388.198
$ws.Range("N136").ClearContents()
